# Refresh the crypto price/volume snapshot as produced by the scheduled
# GitHub Actions job (commit: 'Updated cryptos list ... with GitHub Actions').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.116.46"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.899.43"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'325.33"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4611"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "'0.07885"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'0.9907"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'21.85"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "'1.886.36"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "'7.059"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'5.768"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'0.06990"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "'87.95"
$ws.Range("D17").Value = "'1.002"
$ws.Range("D18").Value = "'0.000009975"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'17.07"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'29.135.89"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'2.111.26"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").Value = "'155.96"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'19.45"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'5.909"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").Value = "'118.34"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("D31").Value = "'0.09335"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'0.8989"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "'5.254"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "'1.325"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'3.157"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").Value = "'0.05790"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'1.175"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'7.723"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").Value = "'0.5693"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "'0.1791"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "'9.715"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "'11.95"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.233"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5354"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'0.07013"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'1.849"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'2.555"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'113.06"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("E51").Value = "  +0.60%  "
